# Edit script implementing the RESTdesign.xlsx update described by the diff:
# "added possibility to add Players to specific Clubs, updated excel design"
#
# Net effect on the data:
#  - Row 12 (endpoint /clubs/{cid}) used to document GET -> "get list of players in
#    each club" and POST -> "add user to club". These are replaced with a richer
#    GET description and a renamed POST description, and a new "GET/POST" note is
#    added in column J. The row grows taller (28.8 -> 43.2) to fit the new text.
#  - Row 22 gets a "done" status note added in column J.
#  - The active cell selection moves from J16 to J12.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: /clubs/{cid} ---------------------------------------------------
# POST description: "add user to club" -> "add player to club"
$ws.Range("F12").Value2 = "add player to club"

# GET description: "get list of players in each club" -> "get information of club,
# and list of players in club"
$ws.Range("E12").Value2 = "get information of club, and list of players in club"

# PUT / DELETE / PATCH columns keep their existing text, they just slide over to
# make room conceptually (content is unchanged - edit club / delete club / edit
# part of club already occupy G12:I12 and stay as-is).

# New note in column J for this row
$ws.Range("J12").Value2 = "GET/POST"

# Row grows to fit the extra wrapped text
$ws.Rows.Item(12).RowHeight = 43.2

# --- Row 22: /admin/users/{id} ---------------------------------------------
$ws.Range("J22").Value2 = "done"

# --- Update the active selection to match the author's final cursor position ---
$ws.Range("J12").Select()
